$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Issues 3-1-24")
$ws.Range("N37").Value = "Done"
$ws.Range("N38").Value = "Done"
